$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns keep their original
# plain-text representation (values such as "0.06850" or "1.010" must not
# be auto-converted to numbers, which would silently drop trailing zeros).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.693.02"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "1.929.65"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").Value = "338.86"
$ws.Range("E5").Value = "  +4.32%  "

$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "0.4092"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("D9").Value = "0.08139"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").Value = "23.51"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "1.942.32"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "6.055"
$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("D14").Value = "7.243"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").Value = "90.74"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").Value = "0.06850"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "0.00001032"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").Value = "17.73"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "1.010"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").Value = "29.716.85"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").Value = "5.610"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").Value = "11.86"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "2.180"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").Value = "2.144.49"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").Value = "6.587"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "157.14"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("D28").Value = "20.01"
$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").Value = "2.084"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").Value = "120.86"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").Value = "1.008"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "0.09635"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").Value = "5.539"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("D34").Value = "1.404"
$ws.Range("E34").Value = "  +2.86%  "

$ws.Range("D35").Value = "3.541"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").Value = "0.06570"
$ws.Range("E36").Value = "  +7.68%  "

$ws.Range("D37").Value = "0.02274"
$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("D38").Value = "1.200"
$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("D39").Value = "0.5949"
$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").Value = "10.72"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").Value = "7.913"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").Value = "2.479"
$ws.Range("E43").Value = "  +2.92%  "

$ws.Range("D44").Value = "1.241"
$ws.Range("E44").Value = "  -3.06%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.07475"
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.24"
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("D47").Value = "0.5540"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").Value = "1.978"
$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("D49").Value = "116.65"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "2.411"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").Value = "72.23"
$ws.Range("E51").Value = "  +0.18%  "
